$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 977
$ws.Range("F4").Value = 842
$ws.Range("F6").Value = 485
$ws.Range("F7").Value = 755
$ws.Range("F8").Value = 175
$ws.Range("F9").Value = 1358
$ws.Range("F10").Value = 784
$ws.Range("F12").Value = 590
$ws.Range("F13").Value = 124
$ws.Range("F14").Value = 3
$ws.Range("F15").Value = 12
$ws.Range("F16").Value = 199
$ws.Range("F17").Value = 91
$ws.Range("F18").Value = 91
$ws.Range("F19").Value = 1443
$ws.Range("F20").Value = 168
$ws.Range("F22").Value = 446
$ws.Range("F23").Value = 33
$ws.Range("F26").Value = 618
$ws.Range("F27").Value = 3
$ws.Range("F28").Value = 182
$ws.Range("F31").Value = 1264
$ws.Range("F32").Value = 64

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 353
$ws.Range("F7").Value = 2
$ws.Range("F8").Value = 267

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 353
$ws.Range("F4").Value = 977
$ws.Range("F6").Value = 842
$ws.Range("F8").Value = 485
$ws.Range("F9").Value = 485
$ws.Range("F10").Value = 755
$ws.Range("F11").Value = 175
$ws.Range("F12").Value = 1358
$ws.Range("F13").Value = 784
$ws.Range("F17").Value = 590
$ws.Range("F19").Value = 124
$ws.Range("F20").Value = 3
$ws.Range("F21").Value = 12
$ws.Range("F22").Value = 199
$ws.Range("F23").Value = 91
$ws.Range("F24").Value = 91
$ws.Range("F25").Value = 1443
$ws.Range("F27").Value = 168
$ws.Range("F29").Value = 446
$ws.Range("F30").Value = 33
$ws.Range("F32").Value = 2
$ws.Range("F34").Value = 267
$ws.Range("F36").Value = 618
$ws.Range("F41").Value = 3
$ws.Range("F42").Value = 182
$ws.Range("F45").Value = 1264
$ws.Range("F46").Value = 64
